$d = $word.ActiveDocument

$replacements = @(
    @("154×8=", "735×3="),
    @("386×7=", "205×9="),
    @("254×7=", "585×2="),
    @("601×9=", "403×2="),
    @("825×5=", "867×9="),
    @("169×4=", "229×2="),
    @("920×6=", "124×8="),
    @("874×2=", "589×9="),
    @("987×2=", "744×2="),
    @("200×6=", "415×7="),
    @("636×4=", "167×4="),
    @("266×7=", "655×3="),
    @("323×3=", "445×8="),
    @("400×2=", "468×4="),
    @("119×4=", "223×7="),
    @("169×2=", "603×2="),
    @("312×4=", "809×9="),
    @("491×6=", "808×8="),
    @("598×4=", "251×3="),
    @("338×6=", "233×5="),
    @("878×4=", "325×8="),
    @("998×9=", "210×2="),
    @("621×5=", "526×8="),
    @("588×4=", "330×9="),
    @("941×6=", "284×2=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
